$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2
$ws.Range("C2").Value = 0.42
$ws.Range("B3").Value = 0.2
$ws.Range("C3").Value = 0.21
$ws.Range("B4").Value = 0.2
$ws.Range("C4").Value = 0.17
$ws.Range("B5").Value = 0.2
$ws.Range("C5").Value = 0.12
$ws.Range("B6").Value = 0.2
$ws.Range("C6").Value = 0.08

$ws.Range("C7").Select()
